$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Save" in H1, reusing the same formatting as the existing header row (G1)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add the new "Save" column values (0) for rows 2-5, reusing formatting of the existing numeric column (G2:G5)
$ws.Range("G2").Copy($ws.Range("H2"))
$ws.Range("H2").Value = 0

$ws.Range("G3").Copy($ws.Range("H3"))
$ws.Range("H3").Value = 0

$ws.Range("G4").Copy($ws.Range("H4"))
$ws.Range("H4").Value = 0

$ws.Range("G5").Copy($ws.Range("H5"))
$ws.Range("H5").Value = 0
